# The workbook contains 28 worksheets (named "41" down to "14"), each with a
# single large text block in cell B2 holding the printed `statsmodels` OLS
# regression summary for one step of a backward-elimination run. Every one of
# those text blocks embeds the run's timestamp ("Date: ..." / "Time: ...").
# This edit re-runs/re-saves the notebook a day later, so every sheet's
# timestamp moves from Wed 01 Jan 2020 23:19:09/23:19:10 to
# Thu 02 Jan 2020 20:49:03 (all 28 occurrences collapse onto the same new
# date/time because the whole report was regenerated in one batch).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value()

    if ($text -ne $null) {
        $updated = $text.Replace("Wed, 01 Jan 2020", "Thu, 02 Jan 2020")
        $updated = $updated.Replace("23:19:09", "20:49:03")
        $updated = $updated.Replace("23:19:10", "20:49:03")

        if ($updated -ne $text) {
            $cell.Value = $updated
        }
    }
}
